$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 2).Value = 0.2109475628992479
$ws.Cells.Item(1, 3).Value = 0.08490505198392986
$ws.Cells.Item(1, 4).Value = 0.7890524371007521
$ws.Cells.Item(1, 5).Value = 0.1944387322314474

$ws.Cells.Item(2, 2).Value = 0.1471589290444466
$ws.Cells.Item(2, 3).Value = 0.03137031465275244
$ws.Cells.Item(2, 4).Value = 0.8528410709555534
$ws.Cells.Item(2, 5).Value = 0.1426829209196241

$ws.Cells.Item(3, 2).Value = 0.916350101644176
$ws.Cells.Item(3, 3).Value = 0.04528999598468715
$ws.Cells.Item(3, 4).Value = 0.083649898355824
$ws.Cells.Item(3, 5).Value = 0.8766467728230319

$ws.Cells.Item(4, 2).Value = 0.0002450104830821196
$ws.Cells.Item(4, 3).Value = 0.0007696321663228415
$ws.Cells.Item(4, 4).Value = 0.9997549895169179
$ws.Cells.Item(4, 5).Value = 0.0002448220601496031

$ws.Cells.Item(5, 2).Value = 0.1204399631601514
$ws.Cells.Item(5, 3).Value = 0.07714737285349238
$ws.Cells.Item(5, 4).Value = 0.8795600368398485
$ws.Cells.Item(5, 5).Value = 0.1118138206484147

$ws.Cells.Item(6, 2).Value = 0.9284232586389959
$ws.Cells.Item(6, 3).Value = 0.1248679750926145
$ws.Cells.Item(6, 4).Value = 0.07157674136100409
$ws.Cells.Item(6, 5).Value = 0.8253619795359144

$ws.Cells.Item(7, 2).Value = 0.8149335533624349
$ws.Cells.Item(7, 3).Value = 0.09273143389216423
$ws.Cells.Item(7, 4).Value = 0.1850664466375651
$ws.Cells.Item(7, 5).Value = 0.7457766181940515

$ws.Cells.Item(8, 2).Value = 0.8337670019782656
$ws.Cells.Item(8, 3).Value = 0.08019858609828262
$ws.Cells.Item(8, 4).Value = 0.1662329980217344
$ws.Cells.Item(8, 5).Value = 0.7718645559330557

$ws.Cells.Item(9, 2).Value = 0.04089765958486836
$ws.Cells.Item(9, 3).Value = 0.02495737852755657
$ws.Cells.Item(9, 4).Value = 0.9591023404151316
$ws.Cells.Item(9, 5).Value = 0.03990181488680195

$ws.Cells.Item(10, 2).Value = 0.6919948240070778
$ws.Cells.Item(10, 3).Value = 0.09746219560244221
$ws.Cells.Item(10, 4).Value = 0.3080051759929222
$ws.Cells.Item(10, 5).Value = 0.630540921391113

$ws.Cells.Item(11, 2).Value = 0.05337750321202273
$ws.Cells.Item(11, 3).Value = 0.003250118200348808
$ws.Cells.Item(11, 4).Value = 0.9466224967879773
$ws.Cells.Item(11, 5).Value = 0.0532045820316198

$ws.Cells.Item(12, 2).Value = 0.1294709107830068
$ws.Cells.Item(12, 3).Value = 0.02574808629088379
$ws.Cells.Item(12, 4).Value = 0.8705290892169932
$ws.Cells.Item(12, 5).Value = 0.1262209625476124

$ws.Cells.Item(13, 2).Value = 0.002815376490472596
$ws.Cells.Item(13, 3).Value = 0.007787633694954902
$ws.Cells.Item(13, 4).Value = 0.9971846235095274
$ws.Cells.Item(13, 5).Value = 0.002793620795038229

$ws.Cells.Item(14, 2).Value = 0.27242409839717
$ws.Cells.Item(14, 3).Value = 0.004044343969174593
$ws.Cells.Item(14, 4).Value = 0.72757590160283
$ws.Cells.Item(14, 5).Value = 0.2713267596530913

$ws.Cells.Item(15, 2).Value = 0.2981947965745748
$ws.Cells.Item(15, 3).Value = 0.1040242579877622
$ws.Cells.Item(15, 4).Value = 0.7018052034254252
$ws.Cells.Item(15, 5).Value = 0.2700980475900741

$ws.Cells.Item(16, 2).Value = 0.584195185130421
$ws.Cells.Item(16, 3).Value = 0.07555270781649608
$ws.Cells.Item(16, 4).Value = 0.415804814869579
$ws.Cells.Item(16, 5).Value = 0.5431581185048652

$ws.Cells.Item(17, 2).Value = 0.01198618216978764
$ws.Cells.Item(17, 3).Value = 0.01606460071924757
$ws.Cells.Item(17, 4).Value = 0.9880138178302124
$ws.Cells.Item(17, 5).Value = 0.01179667332303764

$ws.Cells.Item(18, 2).Value = 0.9274906325732896
$ws.Cells.Item(18, 3).Value = 0.04490490510891151
$ws.Cells.Item(18, 4).Value = 0.07250936742671049
$ws.Cells.Item(18, 5).Value = 0.8876316189525555
